$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed rows 4 and 5 with the same formatting (styles) as row 3, then
# overwrite the values, reproducing the two new review rows.
$ws.Range("A3:G3").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)
$ws.Range("A5:G5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A4").Value = "com.hamxa.shaynachim"
$ws.Range("B4").Value = "bitcoin"
$ws.Range("C4").Value = "erlichyotem@gmail.com "
$ws.Range("D4").Value = "rozend80@gmail.com"
$ws.Range("E4").Value = "27/5/2019 15:59"
$ws.Range("F4").Value = "Things happen fast. Good app"
$ws.Range("G4").Value = "yes"

$ws.Range("A5").Value = "com.hamxa.shaynachim"
$ws.Range("B5").Value = "bitcoin"
$ws.Range("C5").Value = "amramg25@gmail.com "
$ws.Range("D5").Value = "erlichyotem@gmail.com "
$ws.Range("E5").Value = "27/5/2019 15:59"
$ws.Range("F5").Value = "Bitcoin is a fast and furious technology. Great guide"
$ws.Range("G5").Value = "yes"

$ws.Range("C4:D5").Select()
